$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fonts / styles -------------------------------------------------------
# Register a new font (Courier New, black, regular) for the "code echo"
# line. Going through a transient named style avoids the engine's
# Font.Name-setter quirk (it otherwise clones a spare, unused font record
# for the cell's current font before building the new one). Deleting the
# named style afterwards collapses the style back down to a plain cellXf
# that just carries applyFont, leaving fonts/cellXfs with exactly one new
# entry each - matching the shape of the diff (fonts 2->3, cellXfs 2->3).
$codeStyle = $wb.Styles.Add("KnitxlCodeFont")
$codeStyle.Font.Name = "Courier New"

# --- Cell content -----------------------------------------------------
# Previously: A1 held the raw #NUM! error, A3 held the warning message
# (orange/bold), A4 held a leftover empty shared string.
# Now: A1 echoes the evaluated code "log(-1)" in Courier New, A3 keeps the
# warning message/style, and the #NUM! error result moves down to A5 (the
# blank A4 row disappears).
$ws.Range("A1").Style = "KnitxlCodeFont"
$ws.Range("A1").Value = "log(-1)"

$ws.Range("A4").Clear()
$ws.Range("A5").Value = "#NUM!"

# Drop the transient named style now that A1's cellXf carries the font
# directly - keeps cellStyles/cellStyleXfs at their original count (1).
$wb.Styles.Item("KnitxlCodeFont").Delete()
